# Reorganize preset selection rows (8-14) by importance.
# The "CategoryPresetSelection" block rows keep the same Category (A) and
# Control (C) values; only the Option (B), Default (D) and Tooltip (F)
# columns need to be rewritten into their new order.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Menu Mock")

$rows = @(
    @{ Row = 8;  B = "OptionIntensityPreset";         D = "Default";             F = "Intensity profile. Subtle = brief, Default = balanced, Dramatic = stronger, Cinematic = dramatic, Epic = extreme" },
    @{ Row = 9;  B = "OptionTriggerProfile";          D = "All";                 F = "Which triggers are active. Selecting a profile updates the per-trigger toggles." },
    @{ Row = 10; B = "OptionDurationPreset";          D = "Default";             F = "Sets per-trigger duration values." },
    @{ Row = 11; B = "OptionChancePreset";            D = "Off";                 F = "Sets per-trigger chance values. Off means chance is ignored (cooldown only)." },
    @{ Row = 12; B = "OptionCooldownPreset";          D = "Default";             F = "Sets per-trigger cooldown values. Off disables cooldown." },
    @{ Row = 13; B = "OptionTransitionPreset";        D = "Smoothstep";          F = "Sets per-trigger transition curve. Off = instant, Smoothstep = smooth ramp, Linear = constant rate." },
    @{ Row = 14; B = "OptionThirdPersonDistribution"; D = "First Person Only";   F = "Controls how often third-person killcam appears." }
)

foreach ($item in $rows) {
    $r = $item.Row
    $ws.Range("B$r").Value = $item.B
    $ws.Range("D$r").Value = $item.D
    $ws.Range("F$r").Value = $item.F
}
